$d = $word.ActiveDocument

function Replace-Exact($find, $replace) {
    $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false,
        $true, 1, $false, $replace, 2) | Out-Null
}

# 1) "Having " + "around 3+ years..." -> merge into a single run (no text change)
Replace-Exact `
    "Having around 3+ years of hands-on expertise architecting and managing cloud-based services, optimizing deployments using robust CI/CD pipelines using Harness and Jenkin." `
    "Having around 3+ years of hands-on expertise architecting and managing cloud-based services, optimizing deployments using robust CI/CD pipelines using Harness and Jenkin."

# 2) Skills: Cloud Platforms line - add ", Azure cloud " before the final period
Replace-Exact `
    "Cloud Platforms: AWS, ECS, IAM, API Gateway." `
    "Cloud Platforms: AWS, ECS, IAM, API Gateway, Azure cloud ."

# 3) Skills: Agile line - insert "Devops " before "Board."
Replace-Exact `
    "Agile: Jira, Azure Board." `
    "Agile: Jira, Azure Devops Board."

# 4) "Scripting: " + "Python, Bash, YAML, JSON." -> merge into a single run (no text change)
Replace-Exact `
    "Scripting: Python, Bash, YAML, JSON." `
    "Scripting: Python, Bash, YAML, JSON."

# 5) "Result-driven ... and " + "ensure dependable deployment reliability." -> merge (no text change)
Replace-Exact `
    "Result-driven IT professional proficient in Datadog monitoring tools, skilled in designing and executing efficient CI/CD pipelines to streamline software development workflows and ensure dependable deployment reliability." `
    "Result-driven IT professional proficient in Datadog monitoring tools, skilled in designing and executing efficient CI/CD pipelines to streamline software development workflows and ensure dependable deployment reliability."

# 6) " Proficient in J2EE ... cross-country m" + "oney transfers." -> merge (no text change)
Replace-Exact `
    " Proficient in J2EE and Spring Framework, responsible for module development, code optimization, and system performance enhancements. Upgraded legacy code bases to modern standards, aligning functionalities with client requirements. Utilized Jira within a Kanban team for streamlined tracking and management of tasks. Integral in developing the Real-Time Risk Assessment (RTRA) Application, Used tools like Drools for logic implementation of fraudulent transactions in cross-country money transfers." `
    " Proficient in J2EE and Spring Framework, responsible for module development, code optimization, and system performance enhancements. Upgraded legacy code bases to modern standards, aligning functionalities with client requirements. Utilized Jira within a Kanban team for streamlined tracking and management of tasks. Integral in developing the Real-Time Risk Assessment (RTRA) Application, Used tools like Drools for logic implementation of fraudulent transactions in cross-country money transfers."
